# Commit atualizacao tesouro 02/07
# Updates the GRAFICO sheet's cumulative-yield rows (6 and 7), which now
# chain off the previous column's running total instead of recomputing
# from B6/B7 each time, and refreshes the sheet selection/scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GRAFICO")

# --- Row 6 formulas (chain from the previous cumulative column) ---
$ws.Range("D6").Formula  = "=$B$6/E4"
$ws.Range("F6").Formula  = "=D$6+($B$6/G4)"
$ws.Range("H6").Formula  = "=F$6+($B$6/I4)"
$ws.Range("J6").Formula  = "=H$6+($B$6/K4)"
$ws.Range("L6").Formula  = "=J6+($B$6/M4)"
$ws.Range("N6").Formula  = "=L6+($B$6/O4)"
$ws.Range("P6").Formula  = "=N6+($B$6/Q4)"
$ws.Range("R6").Formula  = "=P6+($B$6/S4)"
$ws.Range("T6").Formula  = "=R6+($B$6/U4)"
$ws.Range("V6").Formula  = "=T6+($B$6/W4)"
$ws.Range("X6").Formula  = "=V6+($B$6/Y4)"
$ws.Range("Z6").Formula  = "=X6+($B$6/AA4)"
$ws.Range("AB6").Formula = "=Z6+($B$6/AC4)"
$ws.Range("AD6").Formula = "=AB6+($B$6/AE4)"
$ws.Range("AF6").Formula = "=AD6+($B$6/AG4)"
$ws.Range("AH6").Formula = "=AF6+($B$6/AI4)"
$ws.Range("AJ6").Formula = "=AH6+($B$6/AK4)"
$ws.Range("AL6").Formula = "=AJ6+($B$6/AM4)"
$ws.Range("AN6").Formula = "=AL6+($B$6/AO4)"
$ws.Range("AP6").Formula = "=AN6+($B$6/AQ4)"
$ws.Range("AR6").Formula = "=AP6+($B$6/AS4)"
$ws.Range("AT6").Formula = "=AR6+($B$6/AU4)"

# --- Row 7 formulas ---
$ws.Range("D7").Formula  = "=$B$6/E5"
$ws.Range("F7").Formula  = "=D$6+($B$6/G5)"
$ws.Range("H7").Formula  = "=F$6+($B$6/I5)"
$ws.Range("J7").Formula  = "=H$6+($B$6/K5)"
$ws.Range("L7").Formula  = "=J7+($B$6/M5)"
$ws.Range("N7").Formula  = "=L7+($B$6/O5)"
$ws.Range("P7").Formula  = "=N7+($B$6/Q5)"
$ws.Range("R7").Formula  = "=P7+($B$6/S5)"
$ws.Range("T7").Formula  = "=R7+($B$6/U5)"
$ws.Range("V7").Formula  = "=T7+($B$6/W5)"
$ws.Range("X7").Formula  = "=V7+($B$6/Y5)"
$ws.Range("Z7").Formula  = "=X7+($B$6/AA5)"
$ws.Range("AB7").Formula = "=Z7+($B$6/AC5)"
$ws.Range("AD7").Formula = "=AB7+($B$6/AE5)"
$ws.Range("AF7").Formula = "=AD7+($B$6/AG5)"
$ws.Range("AH7").Formula = "=AF7+($B$6/AI5)"
$ws.Range("AJ7").Formula = "=AH7+($B$6/AK5)"
$ws.Range("AL7").Formula = "=AJ7+($B$6/AM5)"
$ws.Range("AN7").Formula = "=AL7+($B$6/AO5)"
$ws.Range("AP7").Formula = "=AN7+($B$6/AQ5)"
$ws.Range("AR7").Formula = "=AP7+($B$6/AS5)"
$ws.Range("AT7").Formula = "=AR7+($B$6/AU5)"

# D7 used to carry "General" formatting (stale leftover style); it should
# match the rest of the row ("0.00").
$ws.Range("D7").NumberFormat = "0.00"

# --- Sheet view: scroll position + active selection moved ---
$ws.Application.ActiveWindow.ScrollColumn = 42   # column AP
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("AT6").Select()
